$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AC1").Value = "wnb-调节6Hz_20161230_113123_ASIC_EEG"
$ws.Range("AD1").Value = "wnb-调节6Hz_20170110_113300_ASIC_EEG"

$ws.Range("AC2").Value = 0.65916398713826374
$ws.Range("AD2").Value = 0.62135922330097082

$ws.Range("AC3").Value = 0.62732919254658381
$ws.Range("AD3").Value = 0.77815699658703075

$ws.Range("A1:AD3").Select()
